# Updates cryptocurrency Price (D) and Volume(1h) (E) columns for rows 2-51
# A leading apostrophe forces Excel to store the assigned value as literal
# text (matching the source workbook, where these columns are pre-formatted
# strings like "28.875.93" / "0.09330" rather than numbers) without altering
# the cell NumberFormat.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.875.93"
$ws.Range("E2").Value = "'  +1.12%  "

$ws.Range("D3").Value = "'1.877.23"
$ws.Range("E3").Value = "'  -0.85%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "'  -0.72%  "

$ws.Range("D5").Value = "'325.01"
$ws.Range("E5").Value = "'  -0.51%  "

$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "'  -0.61%  "

$ws.Range("D7").Value = "'0.4584"
$ws.Range("E7").Value = "'  -0.19%  "

$ws.Range("D8").Value = "'0.3873"
$ws.Range("E8").Value = "'  +0.27%  "

$ws.Range("D9").Value = "'0.07862"
$ws.Range("E9").Value = "'  -0.09%  "

$ws.Range("D10").Value = "'0.9845"
$ws.Range("E10").Value = "'  -1.78%  "

$ws.Range("D11").Value = "'21.76"
$ws.Range("E11").Value = "'  +0.53%  "

$ws.Range("D12").Value = "'1.908.15"
$ws.Range("E12").Value = "'  +0.87%  "

$ws.Range("D13").Value = "'6.983"
$ws.Range("E13").Value = "'  -1.34%  "

$ws.Range("D14").Value = "'5.643"
$ws.Range("E14").Value = "'  -1.27%  "

$ws.Range("D15").Value = "'0.06962"
$ws.Range("E15").Value = "'  -0.03%  "

$ws.Range("D16").Value = "'88.07"
$ws.Range("E16").Value = "'  +0.74%  "

$ws.Range("D17").Value = "'1.003"
$ws.Range("E17").Value = "'  -0.57%  "

$ws.Range("D18").Value = "'0.000009963"
$ws.Range("E18").Value = "'  -0.86%  "

$ws.Range("D19").Value = "'16.94"
$ws.Range("E19").Value = "'  -1.66%  "

$ws.Range("D20").Value = "'1.003"
$ws.Range("E20").Value = "'  -0.38%  "

$ws.Range("D21").Value = "'28.891.70"
$ws.Range("E21").Value = "'  +1.06%  "

$ws.Range("D22").Value = "'5.249"
$ws.Range("E22").Value = "'  -1.47%  "

$ws.Range("D23").Value = "'10.95"
$ws.Range("E23").Value = "'  -0.56%  "

$ws.Range("D24").Value = "'2.103"
$ws.Range("E24").Value = "'  +2.10%  "

$ws.Range("D25").Value = "'156.02"
$ws.Range("E25").Value = "'  +0.77%  "

$ws.Range("E26").Value = "'  -0.54%  "

$ws.Range("D27").Value = "'6.010"
$ws.Range("E27").Value = "'  +2.63%  "

$ws.Range("D28").Value = "'1.926"
$ws.Range("E28").Value = "'  -1.64%  "

$ws.Range("D29").Value = "'117.21"
$ws.Range("E29").Value = "'  -0.88%  "

$ws.Range("D30").Value = "'0.09330"
$ws.Range("E30").Value = "'  +0.06%  "

$ws.Range("D31").Value = "'0.9014"
$ws.Range("E31").Value = "'  -2.44%  "

$ws.Range("E32").Value = "'  -0.86%  "

$ws.Range("D33").Value = "'1.317"
$ws.Range("E33").Value = "'  -1.45%  "

$ws.Range("D34").Value = "'3.256"
$ws.Range("E34").Value = "'  -0.36%  "

$ws.Range("D35").Value = "'1.183"
$ws.Range("E35").Value = "'  +2.29%  "

$ws.Range("D36").Value = "'0.05752"
$ws.Range("E36").Value = "'  -0.21%  "

$ws.Range("D37").Value = "'0.02065"
$ws.Range("E37").Value = "'  -0.42%  "

$ws.Range("E38").Value = "'  -0.46%  "

$ws.Range("D39").Value = "'7.667"
$ws.Range("E39").Value = "'  -1.66%  "

$ws.Range("D40").Value = "'0.5641"
$ws.Range("E40").Value = "'  -0.51%  "

$ws.Range("D41").Value = "'0.1763"
$ws.Range("E41").Value = "'  -1.40%  "

$ws.Range("D42").Value = "'9.636"
$ws.Range("E42").Value = "'  -1.12%  "

$ws.Range("D43").Value = "'2.255"
$ws.Range("E43").Value = "'  +2.44%  "

$ws.Range("D44").Value = "'11.89"
$ws.Range("E44").Value = "'  +0.89%  "

$ws.Range("D45").Value = "'0.5344"
$ws.Range("E45").Value = "'  -0.22%  "

$ws.Range("D46").Value = "'0.07032"
$ws.Range("E46").Value = "'  -1.82%  "

$ws.Range("D47").Value = "'1.841"
$ws.Range("E47").Value = "'  -0.06%  "

$ws.Range("D48").Value = "'112.88"
$ws.Range("E48").Value = "'  +0.10%  "

$ws.Range("D49").Value = "'2.509"
$ws.Range("E49").Value = "'  +1.46%  "

$ws.Range("D50").Value = "'1.058"
$ws.Range("E50").Value = "'  -5.15%  "

$ws.Range("D51").Value = "'70.62"
$ws.Range("E51").Value = "'  -0.49%  "
